$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the 5 completed replies for DQ1 (row 11, "DQ1 response 5")
$ws.Range("C11").Value = 0.0069444444444444441

# Fill in the completed replies for DQ2 ("DQ2 response 4" / row 16 and
# "DQ2 response 5" / row 17)
$ws.Range("C16").Value = 0.010416666666666666
$ws.Range("C17").Value = 0.010416666666666666

# DQ2 only needed 5 responses, so remove the unused "DQ2 response 6" row
# entirely (row 18) - everything below shifts up by one.
$ws.Rows("18").Delete()

# Move the active selection to C12 (as reflected in the saved file)
$ws.Range("C12").Select()
